$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2: Bitcoin
Set-CellText "D2" "28.866.35"
Set-CellText "E2" "  +7.87%  "

# Row 3: Ethereum
Set-CellText "D3" "1.811.52"
Set-CellText "E3" "  +4.98%  "

# Row 4: TetherUSD
Set-CellText "D4" "0.9990"
Set-CellText "E4" "  +0.17%  "

# Row 5: BNB
Set-CellText "D5" "249.51"
Set-CellText "E5" "  +3.76%  "

# Row 6: USDC
Set-CellText "D6" "0.9995"
Set-CellText "E6" "  +0.18%  "

# Row 7: XRP
Set-CellText "D7" "0.4948"
Set-CellText "E7" "  +1.88%  "

# Row 8: Cardano
Set-CellText "D8" "0.2789"
Set-CellText "E8" "  +7.91%  "

# Row 9: Dogecoin
Set-CellText "D9" "0.06460"

# Row 10: WrappedEther
Set-CellText "D10" "1.811.89"
Set-CellText "E10" "  +4.87%  "

# Row 11: Solana
Set-CellText "E11" "  +5.80%  "

# Row 12: TRON
Set-CellText "D12" "0.07116"
Set-CellText "E12" "  +3.52%  "

# Row 13: Polygon
Set-CellText "E13" "  +6.88%  "

# Row 14: Litecoin
Set-CellText "D14" "84.41"
Set-CellText "E14" "  +9.81%  "

# Row 15: Polkadot
Set-CellText "D15" "4.716"
Set-CellText "E15" "  +5.34%  "

# Row 16: WrappedBTC
Set-CellText "D16" "28.839.33"
Set-CellText "E16" "  +8.57%  "

# Row 17: Dai
Set-CellText "D17" "0.9993"

# Row 18: ShibaInu
Set-CellText "D18" "0.000007420"
Set-CellText "E18" "  +3.74%  "

# Row 19: BinanceUSD
Set-CellText "D19" "0.9988"
Set-CellText "E19" "  +0.14%  "

# Row 20: Avalanche
Set-CellText "D20" "12.27"
Set-CellText "E20" "  +7.33%  "

# Row 21: WrappedliquidstakedEther2.0
Set-CellText "D21" "2.046.66"
Set-CellText "E21" "  +4.90%  "

# Row 22: Uniswap
Set-CellText "D22" "4.599"
Set-CellText "E22" "  +3.91%  "

# Row 23: Cosmos
Set-CellText "D23" "8.984"
Set-CellText "E23" "  +4.71%  "

# Row 24: Chainlink
Set-CellText "D24" "5.359"
Set-CellText "E24" "  +5.35%  "

# Row 25: Monero
Set-CellText "D25" "143.32"
Set-CellText "E25" "  +4.11%  "

# Row 26: BitcoinCash
Set-CellText "D26" "132.22"
Set-CellText "E26" "  +24.83%  "

# Row 27: EthereumClassic
Set-CellText "D27" "16.55"
Set-CellText "E27" "  +8.67%  "

# Row 28: LidoDAOToken
Set-CellText "D28" "1.896"

# Row 29: Toncoin
Set-CellText "D29" "1.395"
Set-CellText "E29" "  +1.73%  "

# Row 30: InternetComputer
Set-CellText "D30" "4.168"
Set-CellText "E30" "  +4.46%  "

# Row 31: Stellar
Set-CellText "D31" "0.08378"
Set-CellText "E31" "  +5.64%  "

# Row 32: Filecoin
Set-CellText "D32" "3.841"
Set-CellText "E32" "  +4.13%  "

# Row 33: Hedera
Set-CellText "D33" "0.04985"
Set-CellText "E33" "  +11.46%  "

# Row 34: ARBITRUM
Set-CellText "D34" "1.094"
Set-CellText "E34" "  +8.74%  "

# Row 35: ImmutableX
Set-CellText "D35" "0.6805"
Set-CellText "E35" "  +9.74%  "

# Row 36: HuobiToken
Set-CellText "D36" "2.706"
Set-CellText "E36" "  +4.27%  "

# Row 37: RenderToken
Set-CellText "D37" "2.278"
Set-CellText "E37" "  +12.62%  "

# Row 38: MXToken
Set-CellText "D38" "2.760"
Set-CellText "E38" "  +13.01%  "

# Row 39: TrustWalletToken
Set-CellText "D39" "0.9589"
Set-CellText "E39" "  +3.71%  "

# Row 40: VeChain
Set-CellText "D40" "0.01597"
Set-CellText "E40" "  +6.82%  "

# Row 41: FraxShare
Set-CellText "D41" "6.048"
Set-CellText "E41" "  +7.05%  "

# Row 42: PaxDollar
Set-CellText "D42" "0.9995"
Set-CellText "E42" "  +0.24%  "

# Row 43: TheSandbox
Set-CellText "D43" "0.4100"
Set-CellText "E43" "  +6.78%  "

# Row 44: Quant
Set-CellText "D44" "99.96"
Set-CellText "E44" "  +0.20%  "

# Row 45: Aptos
Set-CellText "D45" "7.244"
Set-CellText "E45" "  +5.70%  "

# Row 46: Algorand
Set-CellText "D46" "0.1226"
Set-CellText "E46" "  +6.08%  "

# Row 47: Cronos
Set-CellText "E47" "  +2.88%  "

# Row 48: EnergySwap
Set-CellText "D48" "8.142"
Set-CellText "E48" "  +3.43%  "

# Row 49: Elrond
Set-CellText "D49" "31.69"
Set-CellText "E49" "  +5.32%  "

# Row 50: Decentraland
Set-CellText "D50" "0.3642"
Set-CellText "E50" "  +8.41%  "

# Row 51: NEARProtocol
Set-CellText "D51" "1.311"
Set-CellText "E51" "  +6.23%  "
